$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.043795886242208
$ws.Cells.Item(2, 4).Value = 1.04807019863398
$ws.Cells.Item(2, 5).Value = 1.041702456431755
$ws.Cells.Item(2, 6).Value = 1.058069750001876
$ws.Cells.Item(2, 9).Value = 1.038457609427218
$ws.Cells.Item(2, 10).Value = 1.048864414368418
$ws.Cells.Item(2, 11).Value = 1.050831070492827
$ws.Cells.Item(2, 12).Value = 1.04448122347841
$ws.Cells.Item(2, 13).Value = 1.060802991205485
$ws.Cells.Item(2, 14).Value = 1.050353922143251

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.045248276913191
$ws.Cells.Item(3, 4).Value = 1.049151377603628
$ws.Cells.Item(3, 5).Value = 1.042954254642426
$ws.Cells.Item(3, 6).Value = 1.059282513436259
$ws.Cells.Item(3, 9).Value = 1.038765513952434
$ws.Cells.Item(3, 10).Value = 1.049961307204745
$ws.Cells.Item(3, 11).Value = 1.051723523757145
$ws.Cells.Item(3, 12).Value = 1.045542546163096
$ws.Cells.Item(3, 13).Value = 1.061828703637689
$ws.Cells.Item(3, 14).Value = 1.051452372693221

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.046187122098743
$ws.Cells.Item(4, 4).Value = 1.049849976443717
$ws.Cells.Item(4, 5).Value = 1.043763648928317
$ws.Cells.Item(4, 6).Value = 1.060066551454325
$ws.Cells.Item(4, 9).Value = 1.038962965834232
$ws.Cells.Item(4, 10).Value = 1.050669715832677
$ws.Cells.Item(4, 11).Value = 1.052299413263583
$ws.Cells.Item(4, 12).Value = 1.046228137410463
$ws.Cells.Item(4, 13).Value = 1.062491140053672
$ws.Cells.Item(4, 14).Value = 1.052161787342661

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04658159124978
$ws.Cells.Item(5, 4).Value = 1.050143432077041
$ws.Cells.Item(5, 5).Value = 1.044103778008186
$ws.Cells.Item(5, 6).Value = 1.060395996632397
$ws.Cells.Item(5, 9).Value = 1.039045548957054
$ws.Cells.Item(5, 10).Value = 1.050967210535995
$ws.Cells.Item(5, 11).Value = 1.052541140197501
$ws.Cells.Item(5, 12).Value = 1.046516087037975
$ws.Cells.Item(5, 13).Value = 1.062769327881554
$ws.Cells.Item(5, 14).Value = 1.05245970452258

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.046647811541866
$ws.Cells.Item(6, 4).Value = 1.050192690862146
$ws.Cells.Item(6, 5).Value = 1.044160879068674
$ws.Cells.Item(6, 6).Value = 1.060451302385721
$ws.Cells.Item(6, 9).Value = 1.039059390100977
$ws.Cells.Item(6, 10).Value = 1.051017142538805
$ws.Cells.Item(6, 11).Value = 1.052581705192485
$ws.Cells.Item(6, 12).Value = 1.046564419149627
$ws.Cells.Item(6, 13).Value = 1.062816019315068
$ws.Cells.Item(6, 14).Value = 1.052509707434561

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.046192393877519
$ws.Cells.Item(7, 4).Value = 1.049853898532957
$ws.Cells.Item(7, 5).Value = 1.043768194295908
$ws.Cells.Item(7, 6).Value = 1.060070954158682
$ws.Cells.Item(7, 9).Value = 1.038964070985001
$ws.Cells.Item(7, 10).Value = 1.050673692225062
$ws.Cells.Item(7, 11).Value = 1.052302644708385
$ws.Cells.Item(7, 12).Value = 1.046231986075346
$ws.Cells.Item(7, 13).Value = 1.062494858389682
$ws.Cells.Item(7, 14).Value = 1.052165769381979

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.044286928374086
$ws.Cells.Item(8, 4).Value = 1.048435796331541
$ws.Cells.Item(8, 5).Value = 1.042125634766832
$ws.Cells.Item(8, 6).Value = 1.058479756297813
$ws.Cells.Item(8, 9).Value = 1.038562037090171
$ws.Cells.Item(8, 10).Value = 1.049235397260499
$ws.Cells.Item(8, 11).Value = 1.051133009762967
$ws.Cells.Item(8, 12).Value = 1.044840143932087
$ws.Cells.Item(8, 13).Value = 1.061149899995938
$ws.Cells.Item(8, 14).Value = 1.050725431873591

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.040921730484249
$ws.Cells.Item(9, 4).Value = 1.045929145500188
$ws.Cells.Item(9, 5).Value = 1.039226435792704
$ws.Cells.Item(9, 6).Value = 1.055670339356732
$ws.Cells.Item(9, 9).Value = 1.037839893137452
$ws.Cells.Item(9, 10).Value = 1.046690386822725
$ws.Cells.Item(9, 11).Value = 1.049059669733557
$ws.Cells.Item(9, 12).Value = 1.042378524319849
$ws.Cells.Item(9, 13).Value = 1.05877005744075
$ws.Cells.Item(9, 14).Value = 1.048176807229079

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.038672843394641
$ws.Cells.Item(10, 4).Value = 1.044252615305567
$ws.Cells.Item(10, 5).Value = 1.03729015305504
$ws.Cells.Item(10, 6).Value = 1.053793463360658
$ws.Cells.Item(10, 9).Value = 1.037349170395984
$ws.Cells.Item(10, 10).Value = 1.044986368416006
$ws.Cells.Item(10, 11).Value = 1.047668985084953
$ws.Cells.Item(10, 12).Value = 1.040731148388815
$ws.Cells.Item(10, 13).Value = 1.057176669793846
$ws.Cells.Item(10, 14).Value = 1.046470368920768

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037697683170478
$ws.Cells.Item(11, 4).Value = 1.043525325414452
$ws.Cells.Item(11, 5).Value = 1.036450837245209
$ws.Cells.Item(11, 6).Value = 1.052979773550935
$ws.Cells.Item(11, 9).Value = 1.037134460350637
$ws.Cells.Item(11, 10).Value = 1.044246711080027
$ws.Cells.Item(11, 11).Value = 1.047064753705517
$ws.Cells.Item(11, 12).Value = 1.040016270276097
$ws.Cells.Item(11, 13).Value = 1.056485050936561
$ws.Cells.Item(11, 14).Value = 1.045729661186533

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.037335251241452
$ws.Cells.Item(12, 4).Value = 1.04325497200133
$ws.Cells.Item(12, 5).Value = 1.036138938815684
$ws.Cells.Item(12, 6).Value = 1.052677379825227
$ws.Cells.Item(12, 9).Value = 1.037054371807091
$ws.Cells.Item(12, 10).Value = 1.04397169310935
$ws.Cells.Item(12, 11).Value = 1.046840002578618
$ws.Cells.Item(12, 12).Value = 1.039750494693594
$ws.Cells.Item(12, 13).Value = 1.056227897768244
$ws.Cells.Item(12, 14).Value = 1.045454252658792

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.037413003918786
$ws.Cells.Item(13, 4).Value = 1.043312973093887
$ws.Cells.Item(13, 5).Value = 1.036205848466844
$ws.Cells.Item(13, 6).Value = 1.052742251304154
$ws.Cells.Item(13, 9).Value = 1.037071566286665
$ws.Cells.Item(13, 10).Value = 1.044030697980798
$ws.Cells.Item(13, 11).Value = 1.046888226669999
$ws.Cells.Item(13, 12).Value = 1.039807515303896
$ws.Cells.Item(13, 13).Value = 1.056283069605618
$ws.Cells.Item(13, 14).Value = 1.045513341323925

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037667728840877
$ws.Cells.Item(14, 4).Value = 1.043502982134024
$ws.Cells.Item(14, 5).Value = 1.036425058492104
$ws.Cells.Item(14, 6).Value = 1.052954780739702
$ws.Cells.Item(14, 9).Value = 1.037127847063348
$ws.Cells.Item(14, 10).Value = 1.044223983674043
$ws.Cells.Item(14, 11).Value = 1.047046182108307
$ws.Cells.Item(14, 12).Value = 1.039994306066003
$ws.Cells.Item(14, 13).Value = 1.056463799803317
$ws.Cells.Item(14, 14).Value = 1.045706901505026

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037824644702308
$ws.Cells.Item(15, 4).Value = 1.043620025613514
$ws.Cells.Item(15, 5).Value = 1.036560102430467
$ws.Cells.Item(15, 6).Value = 1.053085706721385
$ws.Cells.Item(15, 9).Value = 1.037162478980514
$ws.Cells.Item(15, 10).Value = 1.044343036607334
$ws.Cells.Item(15, 11).Value = 1.047143462117465
$ws.Cells.Item(15, 12).Value = 1.040109362340187
$ws.Cells.Item(15, 13).Value = 1.056575119715332
$ws.Cells.Item(15, 14).Value = 1.045826123507136

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038737531945929
$ws.Cells.Item(16, 4).Value = 1.044300854536116
$ws.Cells.Item(16, 5).Value = 1.037345836443518
$ws.Cells.Item(16, 6).Value = 1.053847444042551
$ws.Cells.Item(16, 9).Value = 1.037363373007251
$ws.Cells.Item(16, 10).Value = 1.045035418562424
$ws.Cells.Item(16, 11).Value = 1.047709042292532
$ws.Cells.Item(16, 12).Value = 1.040778559370476
$ws.Cells.Item(16, 13).Value = 1.057222534701789
$ws.Cells.Item(16, 14).Value = 1.046519488724021

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03930978827277
$ws.Cells.Item(17, 4).Value = 1.044727558585323
$ws.Cells.Item(17, 5).Value = 1.037838464090466
$ws.Cells.Item(17, 6).Value = 1.054324993597161
$ws.Cells.Item(17, 9).Value = 1.037488792032279
$ws.Cells.Item(17, 10).Value = 1.045469244433121
$ws.Cells.Item(17, 11).Value = 1.048063262711398
$ws.Cells.Item(17, 12).Value = 1.041197909821026
$ws.Cells.Item(17, 13).Value = 1.057628190691861
$ws.Cells.Item(17, 14).Value = 1.046953930677214

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.039643443300782
$ws.Cells.Item(18, 4).Value = 1.044976318702931
$ws.Cells.Item(18, 5).Value = 1.03812571968328
$ws.Cells.Item(18, 6).Value = 1.054603444937047
$ws.Cells.Item(18, 9).Value = 1.037561732363994
$ws.Cells.Item(18, 10).Value = 1.045722113796293
$ws.Cells.Item(18, 11).Value = 1.048269675400049
$ws.Cells.Item(18, 12).Value = 1.041442360272428
$ws.Cells.Item(18, 13).Value = 1.057864641966092
$ws.Cells.Item(18, 14).Value = 1.047207159143887

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.039757188797928
$ws.Cells.Item(19, 4).Value = 1.045061117692865
$ws.Cells.Item(19, 5).Value = 1.038223652019774
$ws.Cells.Item(19, 6).Value = 1.054698373625105
$ws.Cells.Item(19, 9).Value = 1.037586566802344
$ws.Cells.Item(19, 10).Value = 1.045808306362332
$ws.Cells.Item(19, 11).Value = 1.048340023308735
$ws.Cells.Item(19, 12).Value = 1.041525686328582
$ws.Cells.Item(19, 13).Value = 1.057945238551941
$ws.Cells.Item(19, 14).Value = 1.047293474113256

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.039248404331451
$ws.Cells.Item(20, 4).Value = 1.04468179066569
$ws.Cells.Item(20, 5).Value = 1.037785618715472
$ws.Cells.Item(20, 6).Value = 1.05427376693176
$ws.Cells.Item(20, 9).Value = 1.03747535795221
$ws.Cells.Item(20, 10).Value = 1.045422717049018
$ws.Cells.Item(20, 11).Value = 1.048025278729908
$ws.Cells.Item(20, 12).Value = 1.041152932990928
$ws.Cells.Item(20, 13).Value = 1.057584684325529
$ws.Cells.Item(20, 14).Value = 1.046907337218889

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037592724664793
$ws.Cells.Item(21, 4).Value = 1.043447034914885
$ws.Cells.Item(21, 5).Value = 1.036360510485857
$ws.Cells.Item(21, 6).Value = 1.052892200379489
$ws.Cells.Item(21, 9).Value = 1.037111283056197
$ws.Cells.Item(21, 10).Value = 1.044167073522158
$ws.Cells.Item(21, 11).Value = 1.046999676831357
$ws.Cells.Item(21, 12).Value = 1.039939307450104
$ws.Cells.Item(21, 13).Value = 1.056410586339854
$ws.Cells.Item(21, 14).Value = 1.045649910534197

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03655049145152
$ws.Cells.Item(22, 4).Value = 1.042669503600139
$ws.Cells.Item(22, 5).Value = 1.035463680668631
$ws.Cells.Item(22, 6).Value = 1.052022666683743
$ws.Cells.Item(22, 9).Value = 1.0368804318374
$ws.Cells.Item(22, 10).Value = 1.043375999520116
$ws.Cells.Item(22, 11).Value = 1.046353029168676
$ws.Cells.Item(22, 12).Value = 1.039174873468775
$ws.Cells.Item(22, 13).Value = 1.055670906224549
$ws.Cells.Item(22, 14).Value = 1.044857713116334

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03710311907066
$ws.Cells.Item(23, 4).Value = 1.04308180195083
$ws.Cells.Item(23, 5).Value = 1.035939185383354
$ws.Cells.Item(23, 6).Value = 1.052483708575918
$ws.Cells.Item(23, 9).Value = 1.037002995123592
$ws.Cells.Item(23, 10).Value = 1.043795516241711
$ws.Cells.Item(23, 11).Value = 1.046696002282787
$ws.Cells.Item(23, 12).Value = 1.039580246693838
$ws.Cells.Item(23, 13).Value = 1.056063166101669
$ws.Cells.Item(23, 14).Value = 1.045277825599793

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.039276141501754
$ws.Cells.Item(24, 4).Value = 1.044702471617891
$ws.Cells.Item(24, 5).Value = 1.037809497531142
$ws.Cells.Item(24, 6).Value = 1.054296914351499
$ws.Cells.Item(24, 9).Value = 1.037481428897546
$ws.Cells.Item(24, 10).Value = 1.045443741308884
$ws.Cells.Item(24, 11).Value = 1.048042442668515
$ws.Cells.Item(24, 12).Value = 1.041173256546919
$ws.Cells.Item(24, 13).Value = 1.05760434347759
$ws.Cells.Item(24, 14).Value = 1.046928391335616

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.041792643708341
$ws.Cells.Item(25, 4).Value = 1.046578116853162
$ws.Cells.Item(25, 5).Value = 1.0399765432819
$ws.Cells.Item(25, 6).Value = 1.056397317832253
$ws.Cells.Item(25, 9).Value = 1.038028217055872
$ws.Cells.Item(25, 10).Value = 1.047349609102093
$ws.Cells.Item(25, 11).Value = 1.049597154423302
$ws.Cells.Item(25, 12).Value = 1.043016004467665
$ws.Cells.Item(25, 13).Value = 1.059386492108155
$ws.Cells.Item(25, 14).Value = 1.048836965679697
